# Adding the changes we made on may 9th
#
# The original sheet held 20 data rows in A2:C21. This update replaces that
# block with 30 data rows (A2:C31): 3 new rows are inserted at the top
# (pushing the former rows 2-21 down to rows 5-24) and 7 more new rows are
# appended after the old data (new rows 25-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 30,3
$data[0,0] = -0.03934990153426188
$data[0,1] = 0.01354811775187661
$data[0,2] = -0.07257660726706182
$data[1,0] = -0.02474731490725574
$data[1,1] = -0.0707512873091867
$data[1,2] = 0.004276057793980506
$data[2,0] = 0.003250675749898542
$data[2,1] = 0.03200497691120401
$data[2,2] = -0.08288132186446855
$data[3,0] = -0.05105815259351689
$data[3,1] = -0.02509637922048548
$data[3,2] = -0.0663297846913338
$data[4,0] = -0.1197514058578581
$data[4,1] = -0.3090105539276482
$data[4,2] = 0.06041020864532091
$data[5,0] = -0.4005748778581615
$data[5,1] = -0.7805985297475537
$data[5,2] = 0.1577123148100716
$data[6,0] = -0.6522004490806942
$data[6,1] = -1.329693669364566
$data[6,2] = 0.1308778794038864
$data[7,0] = -0.6102398293358944
$data[7,1] = -1.220246967815219
$data[7,2] = 0.4809618578070668
$data[8,0] = -0.04692753723689685
$data[8,1] = -0.2711587122508481
$data[8,2] = 1.261436768940515
$data[9,0] = 0.2368920927955998
$data[9,1] = 2.034886604263666
$data[9,2] = 0.6732607796078658
$data[10,0] = -0.7203119397163371
$data[10,1] = 3.954537868499759
$data[10,2] = -2.220546166102086
$data[11,0] = -0.9733701603753284
$data[11,1] = 2.440581185477122
$data[11,2] = -3.891320841653013
$data[12,0] = 0.4988514525549754
$data[12,1] = 0.8867653551555867
$data[12,2] = -1.71078631139937
$data[13,0] = 0.1264272814705268
$data[13,1] = -2.447191684019
$data[13,2] = -0.4622068021978656
$data[14,0] = -0.7295694393771054
$data[14,1] = -0.004014266388761123
$data[14,2] = 0.09374600010258789
$data[15,0] = -0.1179624412740978
$data[15,1] = 1.259269575277969
$data[15,2] = 0.1631955632141657
$data[16,0] = -0.1617193005624273
$data[16,1] = -0.4078179995218952
$data[16,2] = 0.2281727109636559
$data[17,0] = -0.124921940267086
$data[17,1] = 1.249372124671936
$data[17,2] = 1.016479730606079
$data[18,0] = -0.4928955077415405
$data[18,1] = -0.4476696934018742
$data[18,2] = -0.9833766732896989
$data[19,0] = -0.2805471434479678
$data[19,1] = 0.4230750912711692
$data[19,2] = -0.2188279224293611
$data[20,0] = -0.2866266923291341
$data[20,1] = 0.4182899764605935
$data[20,2] = -0.003992439912898826
$data[21,0] = -0.03713915026968551
$data[21,1] = 0.07855436143775912
$data[21,2] = 0.07685266648020034
$data[22,0] = 0.03008511281084441
$data[22,1] = 0.0510363349070151
$data[22,2] = -0.07897615255344478
$data[23,0] = 0.0188495556690863
$data[23,1] = -0.09346238630158486
$data[23,2] = -0.02356194624943445
$data[24,0] = 0.1422297873844697
$data[24,1] = -0.01903863499562003
$data[24,2] = -0.03227404815455254
$data[25,0] = 0.08552113210871153
$data[25,1] = 0.01936588267562919
$data[25,2] = -0.03527019580914864
$data[26,0] = -0.0006981316421711559
$data[26,1] = 0.01939497157166284
$data[26,2] = -0.01994038639324046
$data[27,0] = 0.006530440013323374
$data[27,1] = 0.04903648190555104
$data[27,2] = -0.02060942954960315
$data[28,0] = 0.00994837645529994
$data[28,1] = 0.04445499217226397
$data[28,2] = -0.01362811268440311
$data[29,0] = 0.01760600972920647
$data[29,1] = -0.0257654253925594
$data[29,2] = -0.005214171284543525

$ws.Range("A2:C31").Value2 = $data
